# Updated cryptos list (Price / Volume(1h) refresh) - GitHub Actions style update.
# Rows 14/15 (BitcoinCash <-> WrappedliquidstakedEther2.0) and rows 35/36
# (PEPE <-> Maker) swap coin/link/price/volume data while keeping the
# rank numbers in column A fixed.
# Numeric-looking text in column D is written with a leading apostrophe so
# Excel keeps it as text instead of auto-converting to a number (which would
# otherwise strip meaningful trailing zeros, e.g. "0.220" -> 0.22).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.771.53"
$ws.Range("E2").Value = "  +1.71%  "
$ws.Range("D3").Value = "3.633.21"
$ws.Range("E3").Value = "  +3.77%  "
$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'606.17"
$ws.Range("E5").Value = "  +0.49%  "
$ws.Range("D6").Value = "'199.33"
$ws.Range("E6").Value = "  +2.39%  "
$ws.Range("E7").Value = "  +0.67%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").Value = "'0.220"
$ws.Range("E9").Value = "  +9.33%  "
$ws.Range("E10").Value = "  +0.09%  "
$ws.Range("D11").Value = "'54.05"
$ws.Range("E11").Value = "  +1.61%  "
$ws.Range("E12").Value = "  +2.27%  "
$ws.Range("D13").Value = "'9.58"
$ws.Range("E13").Value = "  +1.22%  "
$ws.Range("B14").Value = "BitcoinCash"
$ws.Range("C14").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D14").Value = "'690.82"
$ws.Range("E14").Value = "  +16.65%  "
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "4.206.50"
$ws.Range("E15").Value = "  +3.52%  "
$ws.Range("D16").Value = "'13.04"
$ws.Range("E16").Value = "  +2.71%  "
$ws.Range("D17").Value = "70.835.32"
$ws.Range("E17").Value = "  +1.50%  "
$ws.Range("D18").Value = "3.659.86"
$ws.Range("E18").Value = "  +4.84%  "
$ws.Range("E19").Value = "  +0.08%  "
$ws.Range("E20").Value = "  +0.35%  "
$ws.Range("E21").Value = "  +1.80%  "
$ws.Range("D22").Value = "'18.85"
$ws.Range("E22").Value = "  +3.48%  "
$ws.Range("D23").Value = "'5.39"
$ws.Range("E23").Value = "  +1.69%  "
$ws.Range("D24").Value = "'105.59"
$ws.Range("E24").Value = "  +4.18%  "
$ws.Range("E25").Value = "  -0.13%  "
$ws.Range("E26").Value = "  -4.30%  "
$ws.Range("D27").Value = "'10.48"
$ws.Range("E27").Value = "  -2.97%  "
$ws.Range("D28").Value = "'9.87"
$ws.Range("E28").Value = "  +4.10%  "
$ws.Range("D29").Value = "'34.25"
$ws.Range("E29").Value = "  +3.40%  "
$ws.Range("D30").Value = "'4.63"
$ws.Range("E30").Value = "  +7.87%  "
$ws.Range("D31").Value = "'7.18"
$ws.Range("E31").Value = "  +2.24%  "
$ws.Range("E32").Value = "  -1.11%  "
$ws.Range("E33").Value = "  +0.55%  "
$ws.Range("D34").Value = "'63.36"
$ws.Range("E34").Value = "  +0.39%  "
$ws.Range("B35").Value = "PEPE"
$ws.Range("C35").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D35").Value = "0.0₃0870"
$ws.Range("E35").Value = "  +6.39%  "
$ws.Range("B36").Value = "Maker"
$ws.Range("C36").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D36").Value = "3.956.14"
$ws.Range("E36").Value = "  +6.39%  "
$ws.Range("D37").Value = "'0.998"
$ws.Range("E37").Value = "  -0.09%  "
$ws.Range("E38").Value = "  -1.62%  "
$ws.Range("D39").Value = "'36.84"
$ws.Range("E39").Value = "  +1.57%  "
$ws.Range("D40").Value = "'504.12"
$ws.Range("E40").Value = "  +3.21%  "
$ws.Range("E41").Value = "  -0.23%  "
$ws.Range("E42").Value = "  -2.92%  "
$ws.Range("D43").Value = "'0.136"
$ws.Range("E43").Value = "  +2.32%  "
$ws.Range("D44").Value = "'3.09"
$ws.Range("E44").Value = "  +10.39%  "
$ws.Range("D45").Value = "'0.0459"
$ws.Range("E45").Value = "  +1.79%  "
$ws.Range("D46").Value = "'3.49"
$ws.Range("E46").Value = "  +6.25%  "
$ws.Range("E47").Value = "  +0.87%  "
$ws.Range("D48").Value = "'8.68"
$ws.Range("E48").Value = "  +3.55%  "
$ws.Range("E49").Value = "  -0.24%  "
$ws.Range("D50").Value = "'0.000248"
$ws.Range("E50").Value = "  +1.47%  "
$ws.Range("E51").Value = "  +1.49%  "
